# audiences-process.pptx — sync non-localizable files
# 1) Refresh the cached "datetimeFigureOut" date field text (slide master +
#    every slide layout) from 01/02/2023 -> 02/08/2023.
# 2) Widen the "Target audience(s) in campaigns" caption textbox on slide 1
#    and extend its wording to mention journeys.

$p = $ppt.ActivePresentation

$oldDate = "01/02/2023"
$newDate = "02/08/2023"

function Update-DateField($shp) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        Update-DateField $shp
    }
}

# Every slide layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            Update-DateField $shp
        }
    }
}

# Slide 1: "Target audience(s) in campaigns" textbox.
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "Target audience(s) in campaigns") {
        $shp.Width = 176.23735
        $shp.TextFrame.TextRange.Text = "Target audience(s) in campaigns and journeys"
    }
}
